$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter headers in the order the strings were first authored (matches the
# shared-string table order produced upstream), then fix up the E/F column
# order to the final displayed layout without re-registering new strings.
$headers = @("destination", "duration", "numberOfTravelers", "specialty", "price", "packageName", "tourLocation", "packageRating", "packageDescription", "contactEmail")

for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Swap columns E1 (price) and F1 (packageName) so the final visible order is
# ... specialty, packageName, price, tourLocation ...
$ws.Cells.Item(1, 5).Value = "packageName"
$ws.Cells.Item(1, 6).Value = "price"

for ($r = 2; $r -le 3; $r++) {
    for ($c = 1; $c -le $headers.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = "test"
    }
}

# Column widths to (approximately) match the authored layout.
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 14.833333333333334
$ws.Columns.Item(6).ColumnWidth = 18
$ws.Columns.Item(7).ColumnWidth = 13.666666666666666
$ws.Columns.Item(8).ColumnWidth = 17.333333333333332
$ws.Columns.Item(9).ColumnWidth = 17.5
$ws.Columns.Item(10).ColumnWidth = 17.166666666666668

# Leave the cursor on I6, matching the saved selection in the workbook.
$ws.Range("I6").Select()
